$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.699.31"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.50%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.879.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.48%  "

$ws.Range("E4").Value = "  +0.49%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "332.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.33%  "

$ws.Range("E6").Value = "  +0.46%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4697"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3941"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.44%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.90"
$ws.Range("D9").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08061"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.94%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.028"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.46%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.14"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.70%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.880.73"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.42%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.979"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.22%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.119"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.41%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.010"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.58%  "

$ws.Range("E17").Value = "  +2.36%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "87.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.44%  "

$ws.Range("E19").Value = "  +1.29%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.64%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.006"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.40%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.538"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.85%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.707.74"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.52%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.79%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.308"
$ws.Range("D25").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.105.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.93%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.79%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.22"
$ws.Range("D28").Style = "Normal"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.101"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.13%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.579"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.27%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "121.77"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.47%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9820"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.12%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09480"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.49%  "

$ws.Range("E34").Value = "  +0.53%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.607"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.56%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.354"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.91%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06138"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.87%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02261"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.70%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.231"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.128"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.79%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5997"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.32%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1898"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.77%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.30"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.46%  "

$ws.Range("E44").Value = "  -0.39%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5718"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.22%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.49%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.947"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.99%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.393"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.39%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06920"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.87%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "114.33"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.01%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.072"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.85%  "
